# Update the multiplication-problem answers in the worksheet table.
# The table has 5 "data" rows (1, 5, 10, 15, 20) each with 5 columns of
# "A x B=C" style text. Because a couple of the original values repeat
# (e.g. "49x70=3430" appears twice but maps to two different new values),
# each replacement is scoped to its specific table cell rather than done
# as a single document-wide Find/Replace.

$mult = [char]0x00D7  # '×' multiplication sign

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Each entry: row, column, old A, old B, old C, new A, new B, new C
$edits = @(
    @(1, 1, "18", "21", "378",  "46", "67", "3082"),
    @(1, 2, "34", "82", "2788", "69", "46", "3174"),
    @(1, 3, "49", "70", "3430", "20", "94", "1880"),
    @(1, 4, "62", "16", "992",  "37", "13", "481"),
    @(1, 5, "45", "36", "1620", "29", "60", "1740"),

    @(5, 1, "87", "98", "8526", "86", "99", "8514"),
    @(5, 2, "63", "84", "5292", "92", "57", "5244"),
    @(5, 3, "95", "18", "1710", "16", "32", "512"),
    @(5, 4, "40", "70", "2800", "46", "52", "2392"),
    @(5, 5, "32", "67", "2144", "55", "21", "1155"),

    @(10, 1, "29", "68", "1972", "64", "81", "5184"),
    @(10, 2, "16", "63", "1008", "21", "68", "1428"),
    @(10, 3, "37", "75", "2775", "34", "23", "782"),
    @(10, 4, "49", "19", "931",  "60", "58", "3480"),
    @(10, 5, "54", "60", "3240", "68", "56", "3808"),

    @(15, 1, "36", "25", "900",  "72", "32", "2304"),
    @(15, 2, "48", "44", "2112", "49", "29", "1421"),
    @(15, 3, "49", "70", "3430", "28", "36", "1008"),
    @(15, 4, "55", "67", "3685", "59", "31", "1829"),
    @(15, 5, "26", "87", "2262", "68", "26", "1768"),

    @(20, 1, "52", "75", "3900", "95", "59", "5605"),
    @(20, 2, "83", "67", "5561", "51", "11", "561"),
    @(20, 3, "92", "76", "6992", "46", "97", "4462"),
    @(20, 4, "53", "71", "3763", "65", "36", "2340"),
    @(20, 5, "27", "54", "1458", "65", "85", "5525")
)

foreach ($edit in $edits) {
    $row = $edit[0]
    $col = $edit[1]
    $searchText  = "{0}{1}{2}={3}" -f $edit[2], $mult, $edit[3], $edit[4]
    $replaceText = "{0}{1}{2}={3}" -f $edit[5], $mult, $edit[6], $edit[7]

    $cellRange = $table.Cell($row, $col).Range
    # Replace:=1 (wdReplaceOne) -- NOT wdReplaceAll(2), which (matching real
    # Word COM semantics) searches/replaces through the *whole* document,
    # ignoring the scoping Range. Since the Find is already scoped to this
    # single cell's Range, "replace one" correctly replaces just the match
    # inside this cell.
    $ok = $cellRange.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 1)
    if (-not $ok) {
        Write-Host "WARNING: replacement failed for row $row col $col : $searchText -> $replaceText"
    }
}

Write-Host "Done."
